# Rewrite M2Doc-style Word fields ( { fldChar begin } { instrText " m: ... " } { fldChar end } )
# into plain text runs containing the literal token text wrapped in curly braces,
# e.g. <w:t>{m: ... }</w:t>, mirroring the new TokenIteratorFieldRewriterSplit
# behaviour (tokens are recognised straight from run text instead of field codes).

$d = $word.ActiveDocument

while ($d.Fields.Count -gt 0) {
    $f = $d.Fields.Item(1)

    # Field instruction text, e.g. " m: 2.myTemplate() " -> "m: 2.myTemplate()"
    $code = $f.Code.Text.Trim()

    # Remember the language formatting carried by the field's code run so the
    # replacement run keeps the same run properties (<w:rPr><w:lang .../></w:rPr>).
    $lang = $f.Code.LanguageID

    # Locate the paragraph that owns this field (the field's Code range sits
    # fully inside exactly one paragraph).
    $codeStart = $f.Code.Start
    $codeEnd = $f.Code.End
    $owner = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $candidate = $d.Paragraphs.Item($i)
        if ($candidate.Range.Start -le $codeStart -and $candidate.Range.End -ge $codeEnd) {
            $owner = $candidate
        }
    }

    # Remove the field (begin/instrText/end runs) entirely, then insert the
    # literal "{code}" text in its place, restoring the original language.
    $f.Delete()
    $owner.Range.InsertBefore("{" + $code + "}")
    $owner.Range.LanguageID = $lang
}
